# Fixed minor issues and updated doc
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: fill in start/end times and the log entry describing the
# report review done together with a team member.
# NOTE: order matters for shared-string placement, so set G70 (and the
# other new text cells) before touching the existing G69 entry.
$ws.Range("D70").Value = 0.98402777777777783
$ws.Range("E70").Value = 0.99097222222222225
$ws.Range("G70").Value = "Reviewed report together with team member for possible issues."

# Row 69: the "Submitted second attempt" note now has a DONE suffix.
$ws.Range("G69").Value = "Submitted second attempt to Canvas. DONE"

# Row 71: fill in start/end times and the revision log entry.
$ws.Range("D71").Value = 0.99097222222222225
$ws.Range("E71").Value = 0.99652777777777779
$ws.Range("G71").Value = "Revised report together with team member for clarity"

# Row 72: fill in start/end times and the final submission log entry.
$ws.Range("D72").Value = 0.99652777777777779
$ws.Range("E72").Value = 0.99930555555555556
$ws.Range("G72").Value = "Last attempt of submission to Canvas. DONE"

# Update the active selection to match the saved view state.
[void]$ws.Range("F69").Select()
